$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data gained one more sampled row for 2026/02/26 (木) that sorts
# ahead of the 2026/12/29 block, so every row from the old row 863 onward
# shifts down by one (old dimension A1:D904 -> new A1:D905).
$ws.Rows.Item(863).Insert()

# Force column A to text first so the date-like string "2026/02/26" is
# stored literally (matching every other date cell in the sheet) instead of
# being auto-parsed into a date serial number by Excel's type inference.
$ws.Cells.Item(863, 1).NumberFormat = "@"
$ws.Cells.Item(863, 1).Value = "2026/02/26"
# Drop back to the sheet's default (unstyled) cell style, same as the other
# data rows, instead of keeping whatever style carried over from the insert.
$ws.Cells.Item(863, 1).Style = "Normal"

$ws.Cells.Item(863, 2).Value = "木"
$ws.Cells.Item(863, 3).Value = 17
$ws.Cells.Item(863, 4).Value = 25
